$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 111-112, pushing the existing rows 111-128 down to 113-130.
$ws.Rows("111:112").Insert()

# New row 111: Early Burlat / Especial
$ws.Cells.Item(111, 1).Value = 8
$ws.Cells.Item(111, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(111, 3).Value = "Coquimbo"
$ws.Cells.Item(111, 4).Value = 44522
$ws.Cells.Item(111, 5).Value = 4
$ws.Cells.Item(111, 6).Value = "Fruta"
$ws.Cells.Item(111, 7).Value = 100103
$ws.Cells.Item(111, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(111, 9).Value = 100103001
$ws.Cells.Item(111, 10).Value = "Cereza"
$ws.Cells.Item(111, 11).Value = "Early Burlat"
$ws.Cells.Item(111, 12).Value = "Especial"
$ws.Cells.Item(111, 13).Value = 200
$ws.Cells.Item(111, 14).Value = 39000
$ws.Cells.Item(111, 15).Value = 40000
$ws.Cells.Item(111, 16).Value = 39500
$ws.Cells.Item(111, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(111, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(111, 19).Value = 2633
$ws.Cells.Item(111, 20).Value = 15

# New row 112: Early Burlat / Primera
$ws.Cells.Item(112, 1).Value = 8
$ws.Cells.Item(112, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(112, 3).Value = "Coquimbo"
$ws.Cells.Item(112, 4).Value = 44522
$ws.Cells.Item(112, 5).Value = 4
$ws.Cells.Item(112, 6).Value = "Fruta"
$ws.Cells.Item(112, 7).Value = 100103
$ws.Cells.Item(112, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(112, 9).Value = 100103001
$ws.Cells.Item(112, 10).Value = "Cereza"
$ws.Cells.Item(112, 11).Value = "Early Burlat"
$ws.Cells.Item(112, 12).Value = "Primera"
$ws.Cells.Item(112, 13).Value = 200
$ws.Cells.Item(112, 14).Value = 36000
$ws.Cells.Item(112, 15).Value = 37000
$ws.Cells.Item(112, 16).Value = 36500
$ws.Cells.Item(112, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(112, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(112, 19).Value = 2433
$ws.Cells.Item(112, 20).Value = 15
